# Apply "improved buttons size, made some basic balancing of units and enemies"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New section title above the existing (friendly) unit table ---
$ws.Range("D3").Value = "Friendly"

# --- Re-balance the friendly unit stats table (rows 5-9, header row 4) ---
# Fishing Boat
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 0.3
$ws.Range("G5").Value = 25
$ws.Range("H5").Value = "-"
$ws.Range("I5").Value = "-"
$ws.Range("J5").Value = 50

# Gunboat
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 0.3
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 0.8
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 100

# Brig
$ws.Range("E7").Value = 20
$ws.Range("F7").Value = 0.3
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 0.8
$ws.Range("I7").Value = 0.5
$ws.Range("J7").Value = 200

# Schooner
$ws.Range("E8").Value = 10
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 10
$ws.Range("H8").Value = 0.7
$ws.Range("I8").Value = 2
$ws.Range("J8").Value = 150

# Warship
$ws.Range("E9").Value = 30
$ws.Range("F9").Value = 0.15
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 0.95
$ws.Range("I9").Value = 1.5
$ws.Range("J9").Value = 250

# --- New "Enemy" table below the friendly table ---
$ws.Range("D13").Value = "Enemy"

$ws.Range("E14").Value = "HP"
$ws.Range("F14").Value = "Speed"
$ws.Range("G14").Value = "Damage"
$ws.Range("H14").Value = "Accuracy"
$ws.Range("I14").Value = "Reload"
$ws.Range("J14").Value = "Price"
$ws.Range("K14").Value = "Damage to ship"

# Enemy Gunboat stats
$ws.Range("D15").Value = "Gunboat"
$ws.Range("E15").Value = 15
$ws.Range("F15").Value = 0.3
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0.8
$ws.Range("I15").Value = 1.1
$ws.Range("K15").Value = 15

# Remaining enemy rows (placeholders, names only so far)
$ws.Range("D16").Value = "Brig"
$ws.Range("D17").Value = "Schooner"
$ws.Range("D18").Value = "Warship"

# --- View/selection tweaks ("improved buttons size") ---
$ws.Range("G8").Select() | Out-Null
